$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "Estudio" (was in X1) to Q1, shifting Congelador..Subposicion (Q1:W1) right by one into R1:X1
$values = @(
    $ws.Range("Q1").Value(),
    $ws.Range("R1").Value(),
    $ws.Range("S1").Value(),
    $ws.Range("T1").Value(),
    $ws.Range("U1").Value(),
    $ws.Range("V1").Value(),
    $ws.Range("W1").Value(),
    $ws.Range("X1").Value()
)

# new order: [Estudio, Congelador, Estante, Posicion rack, Rack, Posicion caja, Caja, Subposicion]
$newValues = @($values[7], $values[0], $values[1], $values[2], $values[3], $values[4], $values[5], $values[6])

$ws.Range("Q1").Value = $newValues[0]
$ws.Range("R1").Value = $newValues[1]
$ws.Range("S1").Value = $newValues[2]
$ws.Range("T1").Value = $newValues[3]
$ws.Range("U1").Value = $newValues[4]
$ws.Range("V1").Value = $newValues[5]
$ws.Range("W1").Value = $newValues[6]
$ws.Range("X1").Value = $newValues[7]

# Update sheet view: scroll position (topLeftCell I1) and selection (N11)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 9
$ws.Range("N11").Select()
